# Change MC of storages
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First table (rows 1-3): rename MC column header, drop the stray "eta"
#     text cells that lived in K1/K2/K3.
$ws.Range("J1").Value = "MC[EUR/MWh]"
$ws.Range("K1").ClearContents()
$ws.Range("K2").ClearContents()
$ws.Range("K3").ClearContents()

# --- Second table (rows 8-10): add new columns for marginal cost in
#     kWh/MWh and the efficiency (eta) that the MC is derived from.
$ws.Range("J8").Value = "MC[EUR/kWh]"
$ws.Range("K8").Value = "MC[EUR/MWh]"
$ws.Range("L8").Value = "eta"

# PumpedStor (row 9) eta = 0.76, BatteryStor (row 10) eta = 0.9
$ws.Range("L9").Value = 0.76
$ws.Range("L10").Value = 0.9

$ws.Range("J9").Formula = "=0.3/L9"
$ws.Range("K9").Formula = "=J9*1000"

$ws.Range("J10").Formula = "=0.3/L10"
$ws.Range("K10").Formula = "=J10*1000"

# Match the style used by the other cells in this section: header cells
# B8/C8/D8/F8/G8/H8 are left-aligned, and the computed annuity cells
# H9/I9/H10/I10 use a 2-decimal number format.
$ws.Range("J8").HorizontalAlignment = $ws.Range("B8").HorizontalAlignment
$ws.Range("K9").NumberFormat = $ws.Range("H9").NumberFormat
$ws.Range("K10").NumberFormat = $ws.Range("H10").NumberFormat

$ws.Range("I14").Select()
